# Update cryptocurrency Price (D) and Volume(1h) (E) columns to refreshed values
# as produced by the GitHub Actions "Updated symbol list" job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.71%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'35.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'13.18%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.152"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.54%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07788"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.91%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.394"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.27%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'8.035"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.03%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.946"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.93%"
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'2.01%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.09819"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'10.75%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.1804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.13%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08665"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.27%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03321"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.63%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.09912"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.35%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001502"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.22%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005699"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.40%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.469"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.47%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.164"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.58%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'1.11%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1294"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.47%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.348"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.67%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2149"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.32%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04570"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.53%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.72%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004448"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.00%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001300"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003699"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'8.92%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'13.37%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.79%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.32%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'6.33%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.007147"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-25.97%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002143"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.93%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.009546"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'12.11%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'25.16%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.001999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.07%"
$ws.Range("E51").Style = "Normal"
